$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Append a new row of data to the "Condicion_Pacientes" table ---------
# This both writes the new cells and grows the table/autoFilter ref
# (A1:F39 -> A1:F40), matching what Excel does when a ListRow is added.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$newRow = $lo.ListRows.Add()

# Copy the date cell's format (built-in date number format) from the row
# above so the new date cell reuses the existing style instead of minting
# a new one.
$ws.Range("A39").Copy($ws.Range("A40"))

$ws.Range("A40").Value2 = 43959
$ws.Range("B40").Value2 = 411
$ws.Range("C40").Value2 = 86
$ws.Range("D40").Value2 = 235
$ws.Range("E40").Value2 = 6
$ws.Range("F40").Value2 = 14

# --- Update the view: scroll position + active selection -----------------
$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1

$ws.Range("D49").Select()
